$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 2 - SNPs in Amplicon
$ws.Range("B2:I2").Value = 10

# Row 3 - SNPs in FW Primer
$ws.Range("B3:I3").Value = 0

# Row 4 - SNPs in RV Primer
$ws.Range("B4:I4").Value = 0

# Row 5 - Max Difference in Melting Temp
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 6
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = "not needed"

# Row 6 - Amplicon Length (from - to)
$ws.Range("B6").Value = "200-500"
$ws.Range("C6").Value = "200-400"
$ws.Range("D6").Value = "200-400"
$ws.Range("E6").Value = "200-500"
$ws.Range("F6").Value = "200-400"
$ws.Range("G6").Value = "200-400"
$ws.Range("H6").Value = "200-500"
$ws.Range("I6").Value = "200-500"

# Row 7 - Min Nr. of GCs per amplicon
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# Row 8 - Min Nr. of CGs per amplicon
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1

# Row 9 - Max Length of monomeric base stretches
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = 7
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 5

# Row 10 - Max length for primer self interaction
$ws.Range("B10:I10").Value = 12

# Row 11 - Primer Length (from - to)
$ws.Range("B11").Value = "18-25"
$ws.Range("C11").Value = "20-32"
$ws.Range("D11").Value = "20-32"
$ws.Range("E11").Value = "18-25"
$ws.Range("F11").Value = "20-32"
$ws.Range("G11").Value = "20-32"
$ws.Range("H11").Value = "18-25"
$ws.Range("I11").Value = "not needed"

# Row 12 - Primer Melting Temp (from - to)
$ws.Range("B12").Value = "50-60"
$ws.Range("C12").Value = "48-60"
$ws.Range("D12").Value = "48-60"
$ws.Range("E12").Value = "50-60"
$ws.Range("F12").Value = "48-60"
$ws.Range("G12").Value = "48-60"
$ws.Range("H12").Value = "50-60"
$ws.Range("I12").Value = "not needed"

# Row 13 - Input Sequence slicing
$ws.Range("B13").Value = 30

# Row 14 - Min C to T conversions FW primer
$ws.Range("C14:H14").Value = 3

# Row 15 - Min G to A conversions RV primer
$ws.Range("C15:H15").Value = 3

# Row 16 - length of one arm in the hp molecule (from - to)
$ws.Range("F16").Value = "50-200"
$ws.Range("G16").Value = "50-200"
$ws.Range("H16").Value = "50-200"

# Selection change to A2
$ws.Range("A2").Select()
